$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Trends Status (sheet1): "Insufficient Data" row count 220 -> 221
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("B8").Value = 221
$wsTrends.Range("C8").Value = 221

# ---------------------------------------------------------------------------
# 2) Priority Status (sheet3): High/Moderate/Low species counts updated
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3) Species qualification (sheet4): rename assessment label + update count
# ---------------------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("A2").Value = "SoIB Assessment"
$wsSpecies.Range("B2").Value = 221

# ---------------------------------------------------------------------------
# 4) "High Priority break-up" sheet becomes two sheets:
#    - Renamed to "Interannual update - High Pri" with new breakdown data
#      (Trend New + IUCN rows, replacing the old single IUCN row)
#    - A brand-new "Major update - High Priority " sheet added right after
#      it, keeping the original break-up content (IUCN / 1 / 100 / 1 / 100)
# ---------------------------------------------------------------------------
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")

# Create the new sheet first (right after the existing one) so it inherits
# the original "High Priority break-up" content before we rewrite it below.
$wsMajor = $wb.Worksheets.Add($null, $wsBreakup)
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"
$wsMajor.Rows.Item(1).Font.Bold = $true
$wsMajor.Rows.Item(1).HorizontalAlignment = -4108

$wsMajor.Range("A2").Value = "IUCN"
$wsMajor.Range("B2").Value = 1
$wsMajor.Range("C2").Value = 100
$wsMajor.Range("D2").Value = 1
$wsMajor.Range("E2").Value = 100

# Now rename the original sheet and rewrite its data with the interannual
# breakdown (a new "Trend New" row plus an updated "IUCN" row).
$wsBreakup.Name = "Interannual update - High Pri"

$wsBreakup.Range("A2").Value = "Trend New"
$wsBreakup.Range("B2").Value = 86
$wsBreakup.Range("C2").Value = 83.5
$wsBreakup.Range("D2").Value = 86
$wsBreakup.Range("E2").Value = 84.3

$wsBreakup.Range("A3").Value = "IUCN"
$wsBreakup.Range("B3").Value = 17
$wsBreakup.Range("C3").Value = 16.5
$wsBreakup.Range("D3").Value = 16
$wsBreakup.Range("E3").Value = 15.7

Write-Output "done"
